$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'long ladies tops',
    'long leg panties',
    'long legging tops for women',
    'long leggings',
    'long leggings for tall women',
    'long leggings for women tall',
    'long length pants',
    'long length yoga pants',
    'long line puff jacket',
    'long running',
    'long running leggings',
    'long running shorts women',
    'long running skirt',
    'long running tights women',
    'long shirts to wear over leggings',
    'long shorts for women below knee',
    'long spanx for women',
    'long string track pants',
    'long tank tops to wear with leggings',
    'long tights',
    'long tights for women',
    'long tops to wear with leggings',
    'long torso compression',
    'long training',
    'long underwear mens skiing',
    'long wear',
    'long workout pants women',
    'loose fitting capris for women',
    'lord normal slides',
    'love clothing brand women',
    'love over distance couple',
    'low rise black tights',
    'low rise capris',
    'low rise compression shorts men',
    'low rise leggings',
    'low rise leggings for women',
    'low rise spandex leggings',
    'low rise tights',
    'low rise tights men',
    'low rise womens leggings',
    'low rise workout capris',
    'low waist brown leggings',
    'low waist tights for women',
    'lower abdomen support',
    'lower back brace xl',
    'lower back compression',
    'lower back compression support',
    'lower back pain yoga',
    'lower back support copper',
    'lower back support gym',
    'lower back support running',
    'lower back support sports',
    'lower back support squat',
    'lower back support women',
    'lower back support workout',
    'lower back support xl',
    'lower belly compression',
    'lucy hatha legging',
    'lucy high rise hatha leggings',
    'lumbar bed support',
    'lumbar compression',
    'lumbar hiking pack',
    'lumbar stabilization exercises',
    'lumbar support black',
    'lumbar support exercise',
    'lumbar support gym',
    'lumbar support hard',
    'lumbar support shirt',
    'lumbar support to wear',
    'lumbar support women',
    'lumbar support workout',
    'lumbar support xl',
    'lumbar warmer',
    'lycra basketball mens',
    'lycra bike shorts women',
    'lycra body',
    'lycra brand',
    'lycra cord',
    'lycra knee sleeve',
    'lycra leggings',
    'lycra leggings for women',
    'lycra leggings women',
    'lycra panties',
    'lycra panties for women',
    'lycra pants',
    'lycra pantyhose for women',
    'lycra running',
    'lycra shorts',
    'lycra shorts for women',
    'lycra shorts women',
    'lycra tights',
    'lycra tights for men',
    'lycra track pants',
    'lycra training',
    'lycra underwear',
    'lycra underwear women',
    'lycra women',
    'lycra workout pants for women',
    'lycras for women gym',
    'lymphedema compression leggings'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
